# "implement the new random item generate item"
#
# The ItemGift table (xl/tables/table1.xml, backed by sheet "ItemGift")
# lists gift-id -> item-drop-list pairs in columns A (Id) and B (Items).
# Rows 12-14 held the old "random generate" gift entries
# (22031101 / 22031102 / 22031103, pointing at the "1;22010001;..." /
# "1;22010101;..." / "1;22010201;..." drop strings). Those are being
# retired, so remove the three rows outright. Excel will shift the rows
# below them (22031201-22031203) up to fill the gap, shrink the table's
# autoFilter/ref range accordingly, and the shared-strings table will
# drop the three now-unused strings on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12:B14").EntireRow.Delete()

# Leave the selection where the next row of data now sits.
$ws.Range("B12").Select()
